# Insert a new record row above row 121 on the active sheet, shifting the
# existing rows 121:140 down to 122:141 (same as a manual "Insert Copied
# Cells" / row-insert in Excel), then populate the new row 121 with the
# latest weekly price observation.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 121:140 down one row -> 122:141
$ws.Rows.Item(121).Insert()

# Populate the newly-inserted row 121 with the new weekly record.
$ws.Range("A121").Value2 = 7
$ws.Range("B121").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C121").Value = "Ñuble"
$ws.Range("D121").Value2 = 44504
$ws.Range("E121").Value2 = 16
$ws.Range("F121").Value2 = 100112017
$ws.Range("G121").Value = "Apio"
$ws.Range("H121").Value = "Americana (o)"
$ws.Range("I121").Value = "Primera"
$ws.Range("J121").Value2 = 120
$ws.Range("K121").Value2 = 8000
$ws.Range("L121").Value2 = 9000
$ws.Range("M121").Value2 = 8500
$ws.Range("N121").Value = "`$/docena de matas"
$ws.Range("O121").Value = "Provincia del Elquí"
$ws.Range("P121").Value2 = 1417
$ws.Range("Q121").Value2 = 6
$ws.Range("R121").Value = "Hortaliza"
